$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$s.Delete()
